$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (after) = old Row 5 values
$ws.Range("D2").Value = 44316
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 20000
$ws.Range("P2").Value = 20000
$ws.Range("S2").Value = 1111

# Row 3 (after) = old Row 2 values
$ws.Range("D3").Value = 44516
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 33000
$ws.Range("O3").Value = 34000
$ws.Range("P3").Value = 33500
$ws.Range("S3").Value = 1861

# Row 4 (after) = old Row 3 values
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 14000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 14500
$ws.Range("S4").Value = 806

# Row 5 (after) = old Row 4 values
$ws.Range("D5").Value = 44280
$ws.Range("L5").Value = "Segunda"
$ws.Range("N5").Value = 12000
$ws.Range("O5").Value = 12000
$ws.Range("P5").Value = 12000
$ws.Range("S5").Value = 667
